# Apply the edits described by the commit:
# "Update TC_ID Excel SCD0017 until SCD0025 and Update TC_ID Solution SCD0006 until SCD0025"
#
# Concretely:
#   1. Rename the worksheet from "SCD0272" to "SCD0017"
#   2. Update cell B2 (TC_ID) from "DGS-287" to "SCD0017-002"
#   3. Select cell B3 (mirrors where the cursor ends up after editing B2 in Excel)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1) Rename the sheet
$ws.Name = "SCD0017"

# 2) Update the TC_ID value in B2
$ws.Range("B2").Value = "SCD0017-002"

# Column B auto-fits to the new, wider text (Excel "best fit" behaviour widened
# the column from 9 to ~13.57 characters once "SCD0017-002" became the longest value)
$ws.Columns.Item(2).ColumnWidth = 12.65

# 3) Move the active selection to B3, as happens after committing an edit in B2
$ws.Range("B3").Select()
